# Fixed up Myxicola and a few more loose ends in second review
#
# Changes applied to the "Materials" sheet:
#   - Remove column A (Taxon_Local_ID / ${iNaturalistTaxonId}) entirely
#   - Remove columns "suborder" / "infraorder" / "superfamily" (both the
#     field-name header row and their ${...} template values) entirely
#   - Rename the stray ${summary.Author} template value to ${summary.authority}

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Materials")

# Fix the typo/rename first, while columns are still in their original
# positions: ${summary.Author} -> ${summary.authority}
# (single-quoted literal so PowerShell doesn't try to expand ${...} as a
# variable reference)
$ws.Range("BB2").Value = '${summary.authority}'

# Delete whole columns from right to left so earlier column letters stay
# valid while we work through them.
$ws.Columns("AT:AT").Delete()
$ws.Columns("AS:AS").Delete()
$ws.Columns("AR:AR").Delete()
$ws.Columns("A:A").Delete()
